$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply cell value updates per the diff. For Price (column D) values that
# are plain decimal numbers (e.g. "244.70"), force a Text number format first
# so Excel keeps them as the exact original text instead of coercing them to
# a floating point number (the source data are textual price strings).

$ws.Range("D2").Value = "36.570.04"
$ws.Range("E2").Value = "  +0.88%  "
$ws.Range("D3").Value = "1.960.36"
$ws.Range("E3").Value = "  -0.46%  "
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "244.70"
$ws.Range("E5").Value = "  +0.61%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.615"
$ws.Range("E6").Value = "  -1.56%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "58.88"
$ws.Range("E7").Value = "  -0.83%  "
$ws.Range("E8").Value = "  -0.07%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.368"
$ws.Range("E9").Value = "  -1.63%  "
$ws.Range("E10").Value = "  -2.15%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0859"
$ws.Range("E11").Value = "  +7.58%  "
$ws.Range("E12").Value = "  +1.04%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "22.21"
$ws.Range("E13").Value = "  -1.10%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.836"
$ws.Range("E14").Value = "  -2.59%  "
$ws.Range("D15").Value = "2.243.76"
$ws.Range("E15").Value = "  -0.73%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "13.73"
$ws.Range("E16").Value = "  -2.39%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.28"
$ws.Range("E17").Value = "  -2.95%  "
$ws.Range("D18").Value = "1.952.11"
$ws.Range("E18").Value = "  -0.92%  "
$ws.Range("D19").Value = "36.465.96"
$ws.Range("E19").Value = "  +0.77%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "70.22"
$ws.Range("E20").Value = "  -1.49%  "
$ws.Range("D21").Value = "0.0₃0872"
$ws.Range("E21").Value = "  +1.57%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "230.77"
$ws.Range("E22").Value = "  -2.56%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.09"
$ws.Range("E23").Value = "  -2.76%  "
$ws.Range("E24").Value = "  +0.04%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.47"
$ws.Range("E25").Value = "  -2.37%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.32"
$ws.Range("E26").Value = "  +1.26%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.37"
$ws.Range("E27").Value = "  -4.75%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.141"
$ws.Range("E28").Value = "  +16.00%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "162.64"
$ws.Range("E29").Value = "  +1.46%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "19.64"
$ws.Range("E30").Value = "  -1.24%  "
$ws.Range("E31").Value = "  -0.78%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.19"
$ws.Range("E32").Value = "  +4.08%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.75"
$ws.Range("E33").Value = "  -2.80%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0638"
$ws.Range("E34").Value = "  +2.75%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.32"
$ws.Range("E35").Value = "  -2.00%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.43"
$ws.Range("E36").Value = "  +2.10%  "
$ws.Range("E37").Value = "  -0.17%  "
$ws.Range("B38").Value = "LidoDAOToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.19"
$ws.Range("E38").Value = "  -4.48%  "
$ws.Range("B39").Value = "WEMIXToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.77"
$ws.Range("E39").Value = "  -2.42%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.07"
$ws.Range("E40").Value = "  -0.48%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0994"
$ws.Range("E41").Value = "  +0.50%  "
$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.19"
$ws.Range("E42").Value = "  -3.23%  "
$ws.Range("B43").Value = "HuobiToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.87"
$ws.Range("E43").Value = "  +0.24%  "
$ws.Range("E44").Value = "  -1.16%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "16.26"
$ws.Range("E45").Value = "  +1.53%  "
$ws.Range("D46").Value = "1.365.28"
$ws.Range("E46").Value = "  +1.77%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.05"
$ws.Range("E47").Value = "  -3.97%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "88.77"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.27"
$ws.Range("E49").Value = "  -3.56%  "
$ws.Range("E50").Value = "  -0.27%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "46.14"
$ws.Range("E51").Value = "  +4.14%  "
